# The workbook used to have two "label" rows at the very top of the
# "Warmup Plan" sheet (row 1 = "Properties", row 2 = "Value") sitting above
# the real header row ("Phase", "Run", "Gmail", ...). This commit removes
# those two now-unused label rows, which shifts every row below them up by
# two and drops the two now out-of-range filler rows at the bottom of the
# sheet (the used range shrinks from A1:T420 to A1:T418). Deleting the rows
# (rather than clearing them) also naturally removes the two now-orphaned
# shared-string entries ("Properties"/"Value") so every other shared-string
# reference (on both sheets) renumbers down by two automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warmup Plan")

$ws.Rows("1:2").Delete()

# The former "Properties"/"Value" rows are gone, so the header row (now row
# 1) is the top of the sheet; mark the whole row as selected, matching the
# saved selection state in the sheet.
[void]$ws.Range("A1:XFD1").Select()
